$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge the separate "NUME" and "PRENUME" headers into a single "NUME PRENUME"
# header that lives in column A.
$ws.Range("A1").Value = "NUME PRENUME"

# The old PRENUME column (B) is no longer needed - delete it entirely. This
# shifts MEDIA FINALA/COD FALCULTATE/SEX/Email/CNP/Telefon one column to the
# left and drops the now-unused trailing column.
$ws.Columns("B").Delete()

$ws.Range("A1").Select()
